$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 data (Improvement: clear/search button)
$ws.Range("B6").Value = "clear button to clear data for new search / take off delete"
$ws.Range("C6").Value = "Improvement"
$ws.Range("D6").Value = 44611
$ws.Range("D6").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E6").Value = "low"
$ws.Range("F6").Value = "yes"

# Update selection to B7
$ws.Range("B7").Select()
